$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the card's name and remaining details into a single Python-tuple-like
# string and store it in A2, then remove the now-redundant rows 3-8.
$combined = "('Jace Beleren', ['{1}{U}{U}', 'Legendary Planeswalker " + [char]0x2014 + " Jace', '+2: Each player draws a card.', '" + [char]0x2212 + "1: Target player draws a card.', '" + [char]0x2212 + "10: Target player mills twenty cards.', 'Loyalty: 3'])"

$ws.Range("A2").Value = $combined

$ws.Range("A3:A8").EntireRow.Delete()
